$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.036.68"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.07%  "

$ws.Range("D3").Value = "'1.708.48"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.46%  "

$ws.Range("D4").Value = "'1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "'307.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -6.16%  "

$ws.Range("D6").Value = "'1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.05%  "

$ws.Range("D7").Value = "'0.4721"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +5.60%  "

$ws.Range("E8").Value = "  -3.87%  "

$ws.Range("D9").Value = "'42.11"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.15%  "

$ws.Range("D10").Value = "'0.07257"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.36%  "

$ws.Range("D11").Value = "'1.034"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.15%  "

$ws.Range("E12").Value = "  +0.01%  "

$ws.Range("D13").Value = "'19.74"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.84%  "

$ws.Range("D14").Value = "'5.837"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.13%  "

$ws.Range("D15").Value = "'1.705.37"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.70%  "

$ws.Range("D16").Value = "'6.824"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.81%  "

$ws.Range("D17").Value = "'88.89"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.78%  "

$ws.Range("E18").Value = "  -2.49%  "

$ws.Range("D19").Value = "'0.06358"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.15%  "

$ws.Range("D20").Value = "'1.002"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.14%  "

$ws.Range("D21").Value = "'16.44"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.03%  "

$ws.Range("D22").Value = "'5.601"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.00%  "

$ws.Range("D23").Value = "'27.076.55"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.12%  "

$ws.Range("D24").Value = "'10.82"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.04%  "

$ws.Range("D25").Value = "'2.111"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.03%  "

$ws.Range("D26").Value = "'156.96"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.04%  "

$ws.Range("D27").Value = "'19.46"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.60%  "

$ws.Range("D28").Value = "'1.901.31"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.71%  "

$ws.Range("D29").Value = "'2.071"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.45%  "

$ws.Range("D30").Value = "'119.19"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.55%  "

$ws.Range("E31").Value = "  -8.98%  "

$ws.Range("D32").Value = "'0.09119"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.04%  "

$ws.Range("D33").Value = "'3.584"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.48%  "

$ws.Range("D34").Value = "'5.287"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.47%  "

$ws.Range("D35").Value = "'0.02185"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.59%  "

$ws.Range("D36").Value = "'0.05801"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.93%  "

$ws.Range("D37").Value = "'10.97"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -7.54%  "

$ws.Range("D38").Value = "'0.1985"
$ws.Range("D38").Style = "Normal"

$ws.Range("D39").Value = "'1.001"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.09%  "

$ws.Range("E40").Value = "  -4.97%  "

$ws.Range("D41").Value = "'1.395"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.21%  "

$ws.Range("D42").Value = "'0.5858"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -7.42%  "

$ws.Range("D43").Value = "'1.097"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -7.29%  "

$ws.Range("D44").Value = "'7.453"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.69%  "

$ws.Range("D45").Value = "'12.58"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.54%  "

$ws.Range("D46").Value = "'0.5632"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.30%  "

$ws.Range("D47").Value = "'3.552"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.04%  "

$ws.Range("D48").Value = "'117.04"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.61%  "

$ws.Range("D49").Value = "'1.829"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.74%  "

$ws.Range("D50").Value = "'0.06623"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.14%  "

$ws.Range("E51").Value = "  -5.09%  "
